$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("ECO list")

# --- Row 12: add new "file names" comment in D12 (SCH/PCB rev note) ---
$ws2.Range("D12").Value = "Plasma_Gen_Main Board_SCH_V2.0_20180107.sch`r`nPlasma_Gen_Main Board_PCB_V2.0_20180107.pcb"
$ws2.Range("D12").WrapText = $true
$ws2.Rows.Item(12).RowHeight = 33

# --- Row 14: new ECO entry - C26 / CAP NC 추가 ---
$ws2.Range("C14").Value = 43107
$ws2.Range("E14").Value = "C26"
$ws2.Range("F9").Copy($ws2.Range("F14"))
$ws2.Range("F9").Copy($ws2.Range("G14"))
$ws2.Range("H14").Value = "CAP NC 추가"
$ws2.Range("I14").Value = "MAX3232 T_OUT pin에 -5V Pull-up이 필요한 경우`r`nC8를  NC처리하고 C26에 0.1uF 실장"
$ws2.Range("I14").WrapText = $true
$ws2.Range("J14").Value = "PCB 수정"
$ws2.Rows.Item(14).RowHeight = 33

# --- Row 15: new ECO entry - R37, R38 / Noise filtering resistor ---
$ws2.Range("C15").Value = 43107
$ws2.Range("E15").Value = "R37, R38"
$ws2.Range("F9").Copy($ws2.Range("F15"))
$ws2.Range("F9").Copy($ws2.Range("G15"))
$ws2.Range("H15").Value = "0 Ω"
$ws2.Range("I15").Value = "Noise filtering or R/Tx pin swap용 저항"
$ws2.Range("J15").Value = "PCB 수정"

# --- Make "ECO list" the active/visible tab, with the last edit cell selected ---
$ws2.Activate()
$ws2.Range("D20").Select()
